# dic_itrf.xlsx edit: update sheet view / selection, resize columns, shrink row 21 height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: drop the frozen/scrolled topLeftCell and move the active selection.
# (Selecting G24 also clears the previous A19:XFD19 full-row selection.)
$ws.Range("G24").Select() | Out-Null

# --- Column widths (values are Excel "character" ColumnWidth units; the stored
# OOXML <col width> is ColumnWidth + 5/6, so each target below is pre-compensated).
$ws.Columns.Item(4).ColumnWidth  = 10.666666666666666   # -> stored 11.5
$ws.Columns.Item(5).ColumnWidth  = 13.666666666666666   # -> stored 14.5
$ws.Columns.Item(6).ColumnWidth  = 15.666666666666666   # -> stored 16.5
$ws.Columns.Item(7).ColumnWidth  = 17.333333333333332   # -> stored ~18.1640625
$ws.Columns.Item(8).ColumnWidth  = 29.666666666666668   # -> stored 30.5
$ws.Columns.Item(12).ColumnWidth = 66.66666666666667    # -> stored 67.5
$ws.Columns.Item(15).ColumnWidth = 24.0                 # -> stored ~24.83203125

# --- Row 21 height shrinks from 22 to 15.
$ws.Rows.Item(21).RowHeight = 15
